# "added support for quoting output"
#
# The sheet used to start at B2 (row2/colB) and ran through H3. The new
# layout starts at A1 and runs through I2 - i.e. every existing cell shifts
# one row up and one column left - and two brand-new trailing columns (H, I)
# are appended with a header + sample value about the new CSV quoting
# behaviour. The mailto: hyperlink (previously on D3) moves along with its
# cell to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shift the whole block up by one row and left by one column ----------
# Deleting row 1 / column 1 pulls every remaining cell (values, number
# formats, wrap text, row heights, the hyperlink style, ...) up/left by one,
# exactly like the old row2/colB content needs to land on row1/colA.
$ws.Rows.Item(1).Delete()
$ws.Columns.Item(1).Delete()

# --- hyperlink: re-point it at its new (shifted) cell ---------------------
# The stored hyperlink range doesn't auto-track the delete above, so drop it
# and recreate it where the Email value now lives.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:edwin@demo.nl")

# --- two new columns: "String with quotes" / "Semicolon in quotes" -------
# Write the data cell before its header so new shared-string entries land in
# the same order as the target workbook (H2, H1, I1, I2).
$ws.Range("H2").Value = 'a "string" containing quotes'
$ws.Range("H1").Value = "String with quotes"
$ws.Range("I1").Value = "Semicolon in quotes"
$ws.Range("I2").Value = 'easy as "1;2;3"'

# --- cosmetic re-save touch-up: give every populated cell an explicit
# Calibri font (mirrors the font/style table normalization seen in the
# target file), without materializing the still-empty D2 cell.
$ws.Range("A1:I1").Font.Name = "Calibri"
$ws.Range("A2:C2").Font.Name = "Calibri"
$ws.Range("E2:I2").Font.Name = "Calibri"

# --- selection, matching the saved workbook's active cell -----------------
$ws.Range("I2").Select()
